# calcolo_perc.xlsx - "modifica per introduzione del test"
# New test run: updated "Mio Algoritmo" (C) results + timestamps (E) for each
# TSP instance, plus refreshed pass/fail highlighting in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$GREEN = 5287936  # RGB(0,176,80) -> FF00B050, same green used across the sheet

# --- row 2 : ch130.tsp  (unchanged result, only the run timestamp moves) ---
$ws.Range("E2").Value = 1556278594271

# --- row 3 : d198.tsp ---
$ws.Range("C3").Value = 15780
$ws.Range("E3").Value = 1556281085299
$ws.Range("F3").Interior.Color = $GREEN
$ws.Range("F3").Font.Color = $GREEN

# --- row 4 : eil76.tsp (unchanged result) ---
$ws.Range("E4").Value = 1556289023497
$ws.Range("F4").Interior.Color = $GREEN
$ws.Range("F4").Font.ColorIndex = -4105

# --- row 5 : fl1577.tsp ---
$ws.Range("C5").Value = 22648
$ws.Range("E5").Value = 1556198455513
$ws.Range("F5").Interior.Color = $GREEN
$ws.Range("F5").Font.ColorIndex = -4105

# --- row 6 : kroA100.tsp (unchanged result) ---
$ws.Range("E6").Value = 1556295011782
$ws.Range("F6").Interior.Color = $GREEN
$ws.Range("F6").Font.ColorIndex = -4105

# --- row 7 : lin318.tsp (unchanged result) ---
$ws.Range("E7").Value = 1556036354013
$ws.Range("F7").Interior.Color = $GREEN
$ws.Range("F7").Font.ColorIndex = -4105

# --- row 8 : pcb442.tsp ---
$ws.Range("C8").Value = 50923
$ws.Range("E8").Value = 1556227370289
$ws.Range("F8").Interior.Color = $GREEN
$ws.Range("F8").Font.ColorIndex = -4105

# --- row 9 : pr439.tsp ---
$ws.Range("C9").Value = 107217
$ws.Range("E9").Value = 1556214127505
$ws.Range("F9").Interior.Color = $GREEN
$ws.Range("F9").Font.Color = $GREEN

# --- row 10 : rat783.tsp ---
$ws.Range("C10").Value = 9124
$ws.Range("E10").Value = 1556388058827
$ws.Range("F10").Interior.Color = $GREEN
$ws.Range("F10").Font.Color = $GREEN

# --- row 11 : u1060.tsp ---
$ws.Range("C11").Value = 226707
$ws.Range("E11").Value = 1556207871133
$ws.Range("F11").Interior.Color = $GREEN
$ws.Range("F11").Font.Color = $GREEN

# D column (shared formula) and the D12 average recalc automatically.

# Leave the selection where the author ended up after entering the last value.
$ws.Range("F11").Select()

Write-Host "calcolo_perc.xlsx updated with new test run"
